$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-7 down to 5-8
$ws.Rows.Item(4).Insert()

# Fill in the new "Waste details" question row
$ws.Cells.Item(4,1).Value = "Waste details"
$ws.Cells.Item(4,2).Value = "waste composition"
$ws.Cells.Item(4,3).Value = "text"
$ws.Cells.Item(4,5).Value = "no"

# Update the active selection to match the saved workbook state
$ws.Range("E4").Select()
